$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Analysis_Unit sheet: update the ExternalInfo() call to pass a parameter ---
$ws1.Range("I3").Value = "LIB_EWS_BE.ExternalInfo(Param1);"

# --- Analysis_Unit sheet: add new "TeradataSchema" column (J) ---
$ws1.Range("J2").Value = "TeradataSchema"
$ws1.Range("J3").Value = "TEWSA0D"

$hdr = $ws1.Range("J1:J2")
$hdr.Font.Name = "Trebuchet MS"
$hdr.Font.Size = 10
$hdr.Font.Bold = $true
$hdr.Font.Color = 16777215
$hdr.Interior.Color = 2499756
$hdr.WrapText = $true

# --- r AnalysisUnit_Variable sheet: add new EXPOSURE variable row ---
$newRow = 108
$ws2.Range("A" + $newRow).Value = "CREATE/MODIFY"
$ws2.Range("B" + $newRow).Value = "COUNTERPARTY_BIB_EXPOSURE"
$ws2.Range("C" + $newRow).Value = "COUNTERPARTY_BIB_EXPOSURE"
$ws2.Range("E" + $newRow).Value = "COUNTERPARTY_BIB"
$ws2.Range("F" + $newRow).Value = "EXPOSURE"
$ws2.Range("A" + $newRow).Font.Name = "Trebuchet MS"

# --- Active tab / selection bookkeeping ---
$ws1.Activate()
$ws1.Range("E8").Select()
$ws2.Range("D111").Select()
